$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row values - written in this order so the shared-string table comes out
# as Hostname, Impact, Possibilite, Mesure de controle (index 0..3).
$ws.Range("A1").Value = "Hostname"
$ws.Range("B1").Value = "Impact"
$ws.Range("C1").Value = "Possibilité"
$ws.Range("D1").Value = "Mesure de contrôle"

# Build the bold 10pt Arial header font on B1 (explicit black text colour), then
# fan it out to C1:D1 via copy/paste-special so they share one style record
# instead of each cell re-deriving its own (avoids duplicate style entries).
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.Size = 10
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Color = 0

$ws.Range("B1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A1 uses the same bold 10pt Arial font but keeps the default (automatic/theme)
# text colour rather than explicit black.
$ws.Range("A1").Font.Name = "Arial"
$ws.Range("A1").Font.Size = 10
$ws.Range("A1").Font.Bold = $true

# Column widths (autofit-like values from the original authoring session).
$ws.Columns.Item(1).ColumnWidth = 11.8
$ws.Columns.Item(3).ColumnWidth = 10.1
$ws.Columns.Item(4).ColumnWidth = 17.8

# Restore the original cursor/selection position.
$ws.Range("D8").Select()
